$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For every data row (2..63) in column A, the date value moves from the
# 1st of the quarter's first month to the 15th of the following month.
for ($row = 2; $row -le 63; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $oldDate = $cell.Value2
    $d = [DateTime]::FromOADate($oldDate)
    $newDate = $d.AddMonths(1).AddDays(14)
    $cell.Value = $newDate.ToOADate()
}
